$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold/bordered/centered) from an existing header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-40 get the season record values
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 70   # AD
    $ws.Cells.Item($r, 31).Value = 92   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
